# The document contains a single table. Its "اولویت" (priority) column
# (column 5) holds a single-digit number for each requirement row, and the
# last column (column 6) holds the requirement code (R1, R2, ... R45) that
# uniquely identifies the row. Update the priority values for the rows
# whose requirement code is a key in $newPriority, leaving every other
# cell (including unlisted rows) untouched.

$newPriority = @{
    "R2"  = "1"
    "R6"  = "1"
    "R10" = "1"
    "R12" = "2"
    "R13" = "2"
    "R14" = "2"
    "R15" = "2"
    "R16" = "3"
    "R17" = "3"
    "R18" = "2"
    "R24" = "2"
    "R25" = "2"
    "R26" = "2"
    "R27" = "3"
    "R28" = "3"
    "R31" = "2"
    "R32" = "3"
    "R33" = "4"
    "R34" = "3"
    "R35" = "4"
    "R36" = "3"
    "R37" = "4"
    "R38" = "3"
    "R39" = "4"
    "R40" = "4"
    "R41" = "4"
    "R42" = "4"
    "R43" = "4"
    "R44" = "4"
    "R45" = "4"
}

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$rowCount = $t.Rows.Count

for ($i = 1; $i -le $rowCount; $i++) {
    $codeCell = $t.Cell($i, 6)
    $codeRange = $codeCell.Range
    [void]$codeRange.MoveEnd(1, -1)
    $code = $codeRange.Text.Trim()

    if ($newPriority.ContainsKey($code)) {
        $priCell = $t.Cell($i, 5)
        $priRange = $priCell.Range
        [void]$priRange.MoveEnd(1, -1)
        $priRange.Text = $newPriority[$code]
    }
}
